$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.473.56"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.094.57"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5212"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4362"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +15.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("D13").Value = "2.092.25"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.671"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.665"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001121"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06588"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "30.520.81"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").Value = "2.340.00"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.558"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.183"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.638"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.137"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.902"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06818"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.437"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2253"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6877"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.260"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6358"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.39%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.625"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.235"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.56%  "
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("E51").Value = "  -1.88%  "
